$wb = $excel.ActiveWorkbook
$wsDoc = $wb.Worksheets.Item("doc vecs")
$wsCls = $wb.Worksheets.Item("classification")

# --- classification sheet: insert two new columns (start lr / end lr) ---
# Shift existing column D ("error rate") to column F by inserting two columns.
$wsCls.Columns("D:E").Insert()

# Header row
$wsCls.Range("D1").Value = "start lr"
$wsCls.Range("E1").Value = "end lr"

# Column widths for the two new columns
# (the engine quantizes ColumnWidth to 1/6-character pixel steps when serializing
#  to OOXML, so these inputs are chosen to land on the closest achievable stored
#  width to the target 9.7109375 / 7.7109375 "characters")
$wsCls.Columns("D").ColumnWidth = 8.8
$wsCls.Columns("E").ColumnWidth = 6.8

# Row 2 (logreg, pvdm)
$wsCls.Range("D2").Value = 5
$wsCls.Range("E2").Value = 0.1
$wsCls.Range("F2").Value = 0.19988

# Row 3 (logreg, dbow)
$wsCls.Range("D3").Value = 5
$wsCls.Range("E3").Value = 0.1
$wsCls.Range("F3").Value = 0.13303999999999999

# Row 4 (logreg, dbow + pvdm)
$wsCls.Range("D4").Value = 5
$wsCls.Range("E4").Value = 0.1
$wsCls.Range("F4").Value = 0.13156000000000001

# Row 5 (neural net, 50)
$wsCls.Range("D5").Value = 1
$wsCls.Range("E5").Value = 1
$wsCls.Range("F5").Value = 0.13496

# Row 6 (neural net, 25)
$wsCls.Range("D6").Value = 1
$wsCls.Range("E6").Value = 1
$wsCls.Range("F6").Value = 0.13456000000000001

# Row 7 (neural net, 10)
$wsCls.Range("D7").Value = 1
$wsCls.Range("E7").Value = 1
$wsCls.Range("F7").Value = 0.13447999999999999

# --- doc vecs sheet: move the selection cursor ---
$wsDoc.Range("A10").Select()

# Restore "classification" as the active sheet/tab
$wsCls.Activate()
